{"js": "const results = context.document.body.search(\"our three zooplankton paper intro\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nconst r = results.items[0];\nconst ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">our three zooplankton paper </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>intro</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\nr.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\nreturn \"ok\";\n", "ps1": "$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Find.MatchCase = $true\n$rng.Find.Execute(\"our three zooplankton paper intro\") | Out-Null\n$para = $rng.Paragraphs(1).Range\n\"Full para range: [\" + $para.Text + \"] start=\" + $para.Start + \" end=\" + $para.End\n$ret = $para.InsertXML('<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Brief summaries of target taxa knowledge, similar to </w:t></w:r><w:r><w:t xml:space=\"preserve\">our three zooplankton paper </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>intro</w:t></w:r><w:proofErr w:type=\"gramEnd\"/></w:p>')\n\"Return: $ret\"\n"}
